{"js": "// Insert a new \"CEPR\" paragraph (yellow highlight) right after the \"BOE\"\n// paragraph, matching the diff that added a <w:p> with a yellow-highlighted\n// run between the \"BOE\" and \"ECB\" paragraphs.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph whose text is exactly \"BOE\".\nlet boeParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"BOE\") {\n    boeParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!boeParagraph) {\n  throw new Error('Could not find a paragraph containing \"BOE\".');\n}\n\n// Insert the new paragraph right after it, with the \"CEPR\" text.\nconst newParagraph = boeParagraph.insertParagraph(\"CEPR\", \"After\");\n\n// Apply a yellow highlight to the new paragraph's text, matching the\n// <w:highlight w:val=\"yellow\"/> run property from the diff.\nnewParagraph.font.highlightColor = \"Yellow\";\n\nawait context.sync();\n", "ps1": "# Insert a new \"CEPR\" paragraph (yellow highlight) right after the \"BOE\"\n# paragraph, matching the diff that added a <w:p> with a yellow-highlighted\n# run between the \"BOE\" and \"ECB\" paragraphs.\n$d = $word.ActiveDocument\n\n$boeParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd() -eq \"BOE\") {\n        $boeParagraph = $p\n    }\n}\n\nif (-not $boeParagraph) {\n    throw 'Could not find a paragraph containing \"BOE\".'\n}\n\n# Insert a new paragraph break right after \"BOE\".\n$boeParagraph.Range.InsertParagraphAfter()\n\n# The newly-created (still empty) paragraph now immediately follows \"BOE\".\n$newParagraph = $boeParagraph.Next()\n$newParagraph.Range.InsertBefore(\"CEPR\")\n\n# Apply a yellow highlight to the new paragraph's text, matching the\n# <w:highlight w:val=\"yellow\"/> run property from the diff (overrides the\n# green highlight that would otherwise be inherited from \"BOE\").\n$newParagraph.Range.HighlightColorIndex = \"wdYellow\"\n"}
